# Apply updated regression summary statistics to the "gw" and "shortage"
# worksheets, and rename a row label on the "gw" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "gw"
# ---------------------------------------------------------------------
$wsGw = $wb.Worksheets.Item("gw")

$wsGw.Range("A6").Value = "l1.cu through l4.cu [NEW]"

$wsGw.Range("B2").Value = 0.4318908365057706
$wsGw.Range("C2").Value = 0.009435258627475629
$wsGw.Range("D2").Value = 0.03422235126330897

$wsGw.Range("B3").Value = 0.5681091634942295
$wsGw.Range("C3").Value = 0.000730393016995956
$wsGw.Range("D3").Value = 0.01373911621442236

$wsGw.Range("B4").Value = 0.8484629433253095
$wsGw.Range("C4").Value = 0.008875022963934831
$wsGw.Range("D4").Value = 0.00008557573900414006

$wsGw.Range("B5").Value = -0.004875163852982749
$wsGw.Range("C5").Value = 0.9390074241580671
$wsGw.Range("D5").Value = 0.8641295925660726

$wsGw.Range("B6").Value = -3.04795110294076
$wsGw.Range("C6").Value = 0.2238945625623958
$wsGw.Range("D6").Value = 0.4441695308514058

$wsGw.Range("B7").Value = 0.00141953030792704
$wsGw.Range("D7").Value = 0.9842707269838914

$wsGw.Range("B9").Value = 0.6673077202139782

# ---------------------------------------------------------------------
# Sheet "shortage"
# ---------------------------------------------------------------------
$wsShortage = $wb.Worksheets.Item("shortage")

$wsShortage.Range("B2").Value = 0.5403018619407929
$wsShortage.Range("C2").Value = 0.00005253177151606476
$wsShortage.Range("D2").Value = 0.0003659900059823309

$wsShortage.Range("B3").Value = 0.9529121868765908
$wsShortage.Range("C3").Value = 0.01514941267356836
$wsShortage.Range("D3").Value = 0.000004579235022134508

$wsShortage.Range("B4").Value = 3.5981532450205
$wsShortage.Range("C4").Value = 0.00003602213589399844
$wsShortage.Range("D4").Value = 0.000001310069128851181

$wsShortage.Range("B6").Value = 2.072908519707469

$wsShortage.Range("B7").Value = 7.827208655252534

$wsShortage.Range("B9").Value = 0.8355209806286752

$wb.Save()
